$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.273.09'
$ws.Range("E2").Value = '  +2.50%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.648.07'
$ws.Range("E3").Value = '  +1.27%  '

$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '603.84'
$ws.Range("E5").Value = '  +2.16%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '156.88'
$ws.Range("E6").Value = '  +4.45%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.11%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.590'
$ws.Range("E8").Value = '  +0.95%  '

$ws.Range("E9").Value = '  +10.63%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.415'
$ws.Range("E10").Value = '  +6.76%  '

$ws.Range("E11").Value = '  +0.56%  '

$ws.Range("E12").Value = '  +1.91%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '29.49'
$ws.Range("E13").Value = '  +6.28%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000192'
$ws.Range("E14").Value = '  +21.97%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.122.86'
$ws.Range("E15").Value = '  +1.24%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.063.27'
$ws.Range("E16").Value = '  +2.49%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.653.45'
$ws.Range("E17").Value = '  +0.85%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.72'
$ws.Range("E18").Value = '  +4.57%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.94'
$ws.Range("E19").Value = '  +3.31%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '361.08'
$ws.Range("E20").Value = '  +4.25%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.43'
$ws.Range("E21").Value = '  +7.61%  '

$ws.Range("E22").Value = '  -0.06%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '69.37'
$ws.Range("E23").Value = '  +3.35%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.72'
$ws.Range("E24").Value = '  +1.26%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.48'
$ws.Range("E25").Value = '  +2.45%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.66'
$ws.Range("E26").Value = '  -0.43%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.38'
$ws.Range("E27").Value = '  -2.20%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0₃0981'
$ws.Range("E28").Value = '  +13.17%  '

$ws.Range("B29").Value = 'Bittensor'
$ws.Range("C29").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '558.15'
$ws.Range("E29").Value = '  +1.84%  '

$ws.Range("B30").Value = 'Kaspa'
$ws.Range("C30").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.166'
$ws.Range("E30").Value = '  +2.56%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.23'
$ws.Range("E31").Value = '  +9.48%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.999'
$ws.Range("E32").Value = '  -0.01%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.83'
$ws.Range("E33").Value = '  +2.62%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.68'
$ws.Range("E34").Value = '  +5.75%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.42'
$ws.Range("E35").Value = '  +4.25%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.434'
$ws.Range("E36").Value = '  +4.95%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '20.54'
$ws.Range("E37").Value = '  +5.01%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.02'
$ws.Range("E38").Value = '  +3.06%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '162.59'
$ws.Range("E39").Value = '  -1.87%  '

$ws.Range("E40").Value = '  -0.03%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.999'
$ws.Range("E41").Value = '  +0.04%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '42.78'
$ws.Range("E42").Value = '  +7.73%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '167.55'
$ws.Range("E43").Value = '  +1.07%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.21'
$ws.Range("E44").Value = '  +2.74%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0626'
$ws.Range("E45").Value = '  +7.33%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.32'
$ws.Range("E46").Value = '  +7.75%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '23.31'
$ws.Range("E47").Value = '  +0.19%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.657'
$ws.Range("E48").Value = '  +3.90%  '

$ws.Range("E49").Value = '  +5.10%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0985'
$ws.Range("E50").Value = '  +2.50%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '19.71'
$ws.Range("E51").Value = '  +2.40%  '
